# Expansão das análises automáticas
# Adds three new computed columns (L: apoio_medio, M: contribuicoes, N: media_contribuicoes)
# to the summary sheet: the header row plus the 6 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (L1:N1) ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the header formatting already used by A1:K1 (bold/centered/bordered)
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Column L - apoio_medio ---
$ws.Range("L2").Value = 92.24142267722461
$ws.Range("L3").Value = 87.74218988390095
$ws.Range("L4").Value = 89.51886785425282
$ws.Range("L5").Value = 94.08839711265563
$ws.Range("L6").Value = 19.15242792744543
$ws.Range("L7").Value = 24.83848419950433

# --- Column M - contribuicoes ---
$ws.Range("M2").Value = 208605
$ws.Range("M3").Value = 54948
$ws.Range("M4").Value = 174765
$ws.Range("M5").Value = 28881
$ws.Range("M6").Value = 2050
$ws.Range("M7").Value = 158

# --- Column N - media_contribuicoes ---
$ws.Range("N2").Value = 312.2829341317365
$ws.Range("N3").Value = 339.1851851851852
$ws.Range("N4").Value = 140.825946817083
$ws.Range("N5").Value = 203.387323943662
$ws.Range("N6").Value = 14.53900709219858
$ws.Range("N7").Value = 14.36363636363636
